# Applies the Tue Apr 4 22:26:12 UTC 2023 cryptos list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.283.34'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '1.874.41'
$ws.Range("E3").Value = '  +3.59%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("D5").Value = "'311.80"
$ws.Range("E5").Value = '  +1.48%  '
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").Value = "'0.5064"
$ws.Range("E7").Value = '  +1.40%  '
$ws.Range("D8").Value = "'0.3935"
$ws.Range("E8").Value = '  +1.25%  '
$ws.Range("D9").Value = "'0.09641"
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("D10").Value = "'1.144"
$ws.Range("E10").Value = '  +4.14%  '
$ws.Range("D11").Value = "'40.96"
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("D12").Value = "'6.497"
$ws.Range("E12").Value = '  +2.98%  '
$ws.Range("E13").Value = '  +2.33%  '
$ws.Range("D14").Value = '1.884.02'
$ws.Range("E14").Value = '  +4.89%  '
$ws.Range("D15").Value = "'7.434"
$ws.Range("E15").Value = '  +3.04%  '
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = "'0.00001132"
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").Value = "'92.90"
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").Value = "'0.06613"
$ws.Range("E19").Value = '  +0.66%  '
$ws.Range("D20").Value = "'17.65"
$ws.Range("E20").Value = '  +3.17%  '
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("D22").Value = "'6.194"
$ws.Range("E22").Value = '  +4.60%  '
$ws.Range("D23").Value = '28.339.68'
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("D24").Value = "'11.29"
$ws.Range("E24").Value = '  +2.74%  '
$ws.Range("D25").Value = "'2.295"
$ws.Range("E25").Value = '  +3.39%  '
$ws.Range("D26").Value = "'2.570"
$ws.Range("E26").Value = '  +6.29%  '
$ws.Range("D27").Value = '2.095.39'
$ws.Range("E27").Value = '  +4.61%  '
$ws.Range("E28").Value = '  +3.03%  '
$ws.Range("D29").Value = "'158.74"
$ws.Range("E29").Value = '  +1.40%  '
$ws.Range("D30").Value = "'127.68"
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").Value = "'1.070"
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("D33").Value = "'5.633"
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("E34").Value = '  +0.76%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = "'0.06737"
$ws.Range("E35").Value = '  -1.01%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = "'9.510"
$ws.Range("E36").Value = '  +7.02%  '
$ws.Range("D37").Value = "'0.02404"
$ws.Range("E37").Value = '  +4.28%  '
$ws.Range("D38").Value = "'0.2191"
$ws.Range("E38").Value = '  +2.46%  '
$ws.Range("D39").Value = "'11.52"
$ws.Range("E39").Value = '  +1.23%  '
$ws.Range("D40").Value = "'0.6374"
$ws.Range("E40").Value = '  +2.96%  '
$ws.Range("D41").Value = "'5.000"
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("D42").Value = "'1.183"
$ws.Range("E42").Value = '  +3.51%  '
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("D44").Value = "'13.56"
$ws.Range("E44").Value = '  +3.63%  '
$ws.Range("D45").Value = "'0.5997"
$ws.Range("E45").Value = '  +2.08%  '
$ws.Range("D46").Value = "'3.660"
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").Value = "'1.267"
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("D48").Value = "'2.005"
$ws.Range("E48").Value = '  +3.32%  '
$ws.Range("D49").Value = "'124.22"
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("D50").Value = "'1.196"
$ws.Range("E50").Value = '  +1.92%  '
$ws.Range("D51").Value = "'0.06857"
$ws.Range("E51").Value = '  +2.04%  '
